$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.90532904932639
$ws.Range("C2").Value = 10.58136980179555
$ws.Range("D2").Value = 4.729078971557712
$ws.Range("E2").Value = 11.24430392303668
$ws.Range("F2").Value = 23.43303130976396
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 20.96508356542164
$ws.Range("L2").Value = 9.841459993918596
$ws.Range("O2").Value = 20.96543010837515
$ws.Range("B3").Value = 16.20083200353825
$ws.Range("C3").Value = 10.26389082591426
$ws.Range("D3").Value = 4.684323851078094
$ws.Range("E3").Value = 11.29889216155078
$ws.Range("F3").Value = 23.48647448217672
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 21.12187738033679
$ws.Range("L3").Value = 9.807471450699266
$ws.Range("O3").Value = 21.06552567470003
$ws.Range("B4").Value = 15.75268506530987
$ws.Range("C4").Value = 10.06300304160428
$ws.Range("D4").Value = 4.65648703907255
$ws.Range("E4").Value = 11.3345345646507
$ws.Range("F4").Value = 23.52808438401473
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 21.22394490102202
$ws.Range("L4").Value = 9.788220897778785
$ws.Range("O4").Value = 21.13363487981578
$ws.Range("B5").Value = 15.56637990378482
$ws.Range("C5").Value = 9.979733876691393
$ws.Range("D5").Value = 4.645060004190038
$ws.Range("E5").Value = 11.34959383497555
$ws.Range("F5").Value = 23.54724161465516
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 21.26699470959756
$ws.Range("L5").Value = 9.780788611670523
$ws.Range("O5").Value = 21.1630541671847
$ws.Range("B6").Value = 15.53522923935286
$ws.Range("C6").Value = 9.965825018623327
$ws.Range("D6").Value = 4.643157711752322
$ws.Range("E6").Value = 11.3521267207185
$ws.Range("F6").Value = 23.55055524329716
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 21.27423104176459
$ws.Range("L6").Value = 9.779579550104538
$ws.Range("O6").Value = 21.16803949636353
$ws.Range("B7").Value = 15.7501870600402
$ws.Range("C7").Value = 10.06188560973205
$ws.Range("D7").Value = 4.656333258840135
$ws.Range("E7").Value = 11.33473549376834
$ws.Range("F7").Value = 23.52833384922659
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 21.22451958969515
$ws.Range("L7").Value = 9.788118986372355
$ws.Range("O7").Value = 21.13402491090545
$ws.Range("B8").Value = 16.66578540602326
$ws.Range("C8").Value = 10.47319138236834
$ws.Range("D8").Value = 4.713725355616287
$ws.Range("E8").Value = 11.2626851131501
$ws.Range("F8").Value = 23.44962683290585
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 21.0179426817926
$ws.Range("L8").Value = 9.829408220528132
$ws.Range("O8").Value = 20.99855806804045
$ws.Range("B9").Value = 18.32871560163888
$ws.Range("C9").Value = 11.22888965123323
$ws.Range("D9").Value = 4.823159012711181
$ws.Range("E9").Value = 11.13823946873638
$ws.Range("F9").Value = 23.36549886483291
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 20.65887061883354
$ws.Range("L9").Value = 9.922960637898019
$ws.Range("O9").Value = 20.7860239323911
$ws.Range("B10").Value = 19.45960672573077
$ws.Range("C10").Value = 11.7487709134345
$ws.Range("D10").Value = 4.901295433160203
$ws.Range("E10").Value = 11.05705609475013
$ws.Range("F10").Value = 23.34699124613472
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 20.42317635572225
$ws.Range("L10").Value = 9.998994957977905
$ws.Range("O10").Value = 20.66272354443374
$ws.Range("B11").Value = 19.95269727419245
$ws.Range("C11").Value = 11.97683905555691
$ws.Range("D11").Value = 4.936274066616118
$ws.Range("E11").Value = 11.02234437297453
$ws.Range("F11").Value = 23.34804994690534
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 20.32207362663503
$ws.Range("L11").Value = 10.03508375501266
$ws.Range("O11").Value = 20.61385969187142
$ws.Range("B12").Value = 20.13623565170902
$ws.Range("C12").Value = 12.06193542354989
$ws.Range("D12").Value = 4.949431923421529
$ws.Range("E12").Value = 11.00951880747878
$ws.Range("F12").Value = 23.34981750186159
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 20.28466989770161
$ws.Range("L12").Value = 10.04895753269616
$ws.Range("O12").Value = 20.59640288224698
$ws.Range("B13").Value = 20.09685061905299
$ws.Range("C13").Value = 12.04366562360662
$ws.Range("D13").Value = 4.946602148411904
$ws.Range("E13").Value = 11.0122668354231
$ws.Range("F13").Value = 23.34937600292963
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 20.29268620941865
$ws.Range("L13").Value = 10.04596045842639
$ws.Range("O13").Value = 20.60011583799427
$ws.Range("B14").Value = 19.96786149169693
$ws.Range("C14").Value = 11.98386569843761
$ws.Range("D14").Value = 4.937358367026802
$ws.Range("E14").Value = 11.02128281290904
$ws.Range("F14").Value = 23.34816795939204
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 20.31897871569378
$ws.Range("L14").Value = 10.03622104280522
$ws.Range("O14").Value = 20.61240248529738
$ws.Range("B15").Value = 19.88843417511356
$ws.Range("C15").Value = 11.94706981723133
$ws.Range("D15").Value = 4.931684663476476
$ws.Range("E15").Value = 11.02684690206282
$ws.Range("F15").Value = 23.34760605911153
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 20.33519852529252
$ws.Range("L15").Value = 10.03028217919816
$ws.Range("O15").Value = 20.62006497075814
$ws.Range("B16").Value = 19.4269422568432
$ws.Range("C16").Value = 11.73369158361973
$ws.Range("D16").Value = 4.898997557000823
$ws.Range("E16").Value = 11.05936923201921
$ws.Range("F16").Value = 23.3471131333948
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 20.42990690922025
$ws.Range("L16").Value = 9.99666598348993
$ws.Range("O16").Value = 20.66606298994536
$ws.Range("B17").Value = 19.13827556553996
$ws.Range("C17").Value = 11.60059029852836
$ws.Range("D17").Value = 4.878795505108116
$ws.Range("E17").Value = 11.07988890312782
$ws.Range("F17").Value = 23.34924142666565
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 20.48957536882259
$ws.Range("L17").Value = 9.976422217531237
$ws.Range("O17").Value = 20.69613767597863
$ws.Range("B18").Value = 18.97023740928921
$ws.Range("C18").Value = 11.52324357080211
$ws.Range("D18").Value = 4.86712291689512
$ws.Range("E18").Value = 11.09190014139592
$ws.Range("F18").Value = 23.35135754760685
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 20.5244705891447
$ws.Range("L18").Value = 9.964920264680845
$ws.Range("O18").Value = 20.71411554787159
$ws.Range("B19").Value = 18.91300200461406
$ws.Range("C19").Value = 11.49692136033844
$ws.Range("D19").Value = 4.86316189051724
$ws.Range("E19").Value = 11.09600281610157
$ws.Range("F19").Value = 23.35222708863803
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 20.53638428864271
$ws.Range("L19").Value = 9.961050481943767
$ws.Range("O19").Value = 20.72031904589734
$ws.Range("B20").Value = 19.16921298280129
$ws.Range("C20").Value = 11.6148413860413
$ws.Range("D20").Value = 4.88095157003614
$ws.Range("E20").Value = 11.07768293289811
$ws.Range("F20").Value = 23.34892251712627
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 20.48316398159875
$ws.Range("L20").Value = 9.97856259000808
$ws.Range("O20").Value = 20.69286576994489
$ws.Range("B21").Value = 20.00583596207952
$ws.Range("C21").Value = 12.00146521764904
$ws.Range("D21").Value = 4.940075922607199
$ws.Range("E21").Value = 11.01862594443466
$ws.Range("F21").Value = 23.34848567885327
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 20.31123202197094
$ws.Range("L21").Value = 10.03907617288775
$ws.Range("O21").Value = 20.60876512686871
$ws.Range("B22").Value = 20.53401616561882
$ws.Range("C22").Value = 12.24673314597931
$ws.Range("D22").Value = 4.978202524179983
$ws.Range("E22").Value = 10.98188820513257
$ws.Range("F22").Value = 23.35616663698141
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 20.20400482268751
$ws.Range("L22").Value = 10.07983219591619
$ws.Range("O22").Value = 20.55990586128394
$ws.Range("B23").Value = 20.25385145608353
$ws.Range("C23").Value = 12.11652427415695
$ws.Range("D23").Value = 4.957902797976518
$ws.Range("E23").Value = 11.00132571099462
$ws.Range("F23").Value = 23.35133741793375
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 20.26076287480625
$ws.Range("L23").Value = 10.057972236207
$ws.Range("O23").Value = 20.58542183686688
$ws.Range("B24").Value = 19.15523264591534
$ws.Range("C24").Value = 11.60840103487054
$ws.Range("D24").Value = 4.879976993739223
$ws.Range("E24").Value = 11.07867958479492
$ws.Range("F24").Value = 23.34906391621479
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 20.48606072886177
$ws.Range("L24").Value = 9.977594501959382
$ws.Range("O24").Value = 20.69434285758994
$ws.Range("B25").Value = 17.89417387060793
$ws.Range("C25").Value = 11.0303738606398
$ws.Range("D25").Value = 4.793925753461072
$ws.Range("E25").Value = 11.17010434038269
$ws.Range("F25").Value = 23.38068381269725
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 20.75107646186859
$ws.Range("L25").Value = 9.896341693515613
$ws.Range("O25").Value = 20.8377845566254
